$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: seed new shared strings in the exact order they first appear so the
# resulting sharedStrings.xml table matches the target ordering:
#   25 transformation, 26 log10, 27 squared, 28 NA, 29 MeanResponse,
#   30 logit (beta models)
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "transformation"
$ws.Range("F5").Value = "log10"
$ws.Range("F11").Value = "squared"
$ws.Range("F2").Value = "NA"
$ws.Range("D1").Value = "MeanResponse"
$ws.Range("F8").Value = "logit (beta models)"
$ws.Range("E1").Value = "percEstimate"

# ---------------------------------------------------------------------------
# Step 2: column D - the old "percEstimate" formulas are replaced either by
# plain numeric mean-response values (for rows that used a ratio-based
# back-transform) or removed entirely (for rows that used a simple *100 or
# EXP back-transform, which no longer need a helper column).
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 27.287123139999999
$ws.Range("D3").Value = 19.245025900000002
$ws.Range("D4").Value = 1.9782322990000001
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("D7").Value = 0.54293340099999998
$ws.Range("D8").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("D10").Value = 0.51791088799999996
$ws.Range("D11").Value = 0.82625772399999997
$ws.Range("D12").Value = 37.709882800000003
$ws.Range("D13").Value = 0.29114584599999999
$ws.Range("D14").Value = 1.42520523
$ws.Range("D15").ClearContents()
$ws.Range("D16").Value = 27.528243620000001
$ws.Range("D17").ClearContents()
$ws.Range("D18").Value = 9.2795485230000008
$ws.Range("D19").ClearContents()
$ws.Range("D20").Value = 19.19656517
$ws.Range("D21").ClearContents()

# ---------------------------------------------------------------------------
# Step 3: column E now carries the percent-change-per-year calculation
# (corrected transformation logic).
# ---------------------------------------------------------------------------
$ws.Range("E2").Formula = "=(C2/D2)*100"
$ws.Range("E3:E4").Formula = "=(C3/D3)*100"
$ws.Range("E5").Formula = "=(10^C5-1)*100"
$ws.Range("E6").Formula = "=(10^C6-1)*100"
$ws.Range("E7").Formula = "=(C7/D7)*100"
$ws.Range("E8").Formula = "=(EXP(C8)-1)*100"
$ws.Range("E9").Formula = "=(10^C9-1)*100"
$ws.Range("E10").Formula = "=(C10/D10)*100"
$ws.Range("E11").Formula = "=(C11/(D11^2))*100"
$ws.Range("E12:E14").Formula = "=(C12/D12)*100"
$ws.Range("E13").ClearFormats()
$ws.Range("E15").Formula = "=(10^C15-1)*100"
$ws.Range("E16").Formula = "=(C16/D16)*100"
$ws.Range("E17").Formula = "=(10^C17-1)*100"
$ws.Range("E18").Formula = "=(C18/D18)*100"
$ws.Range("E19").Formula = "=(10^C19-1)*100"
$ws.Range("E20").Formula = "=(C20/D20)*100"
$ws.Range("E21").Formula = "=(10^C21-1)*100"

# ---------------------------------------------------------------------------
# Step 4: column F records which back-transformation was used for each row.
# ---------------------------------------------------------------------------
$ws.Range("F3").Value = "NA"
$ws.Range("F4").Value = "NA"
$ws.Range("F6").Value = "log10"
$ws.Range("F7").Value = "NA"
$ws.Range("F9").Value = "log10"
$ws.Range("F10").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("F13").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("F15").Value = "log10"
$ws.Range("F16").Value = "NA"
$ws.Range("F17").Value = "log10"
$ws.Range("F18").Value = "NA"
$ws.Range("F19").Value = "log10"
$ws.Range("F20").Value = "NA"
$ws.Range("F21").Value = "log10"

# ---------------------------------------------------------------------------
# Step 5: cosmetic worksheet view changes - zoom, selection and column widths.
# ---------------------------------------------------------------------------
$ws.Range("A:A").ColumnWidth = 4.6640625
$ws.Range("B:B").ColumnWidth = 17.25
$excel.ActiveWindow.Zoom = 86
$ws.Range("G14").Select()
